# Auto update: 2025-12-05 12:21:03
# Update MACRO_SCORE (column N) values on rows 2-6 of Sheet1 from
# 54.86376272656823 to 54.84087454262382

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 54.84087454262382

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 14).Value = $newValue
}
